$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write "123" into A2 as literal text (matching the original string cell)
# without disturbing A2's own style/number format. A direct .Value = "123"
# would be auto-coerced to a number by Excel's input parsing, so stage the
# text in a scratch cell formatted as Text, then paste-special just the
# value (xlPasteValues = -4163) onto A2 -- that carries the text type over
# but leaves A2's existing formatting alone. Finally wipe the scratch cell.
$scratch = $ws.Range("ZZ500")
$scratch.NumberFormat = "@"
$scratch.Value = "123"
$scratch.Copy()
$ws.Cells.Item(2, 1).PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = $false
$ws.Range("B2").Value = 0.8620919585227966
$ws.Range("F2").Value = 0.8718597888946533
$ws.Range("G2").Value = 97.02018737792969
$ws.Range("H2").Value = 14.78398036956787
$ws.Range("I2").Value = 12.06411457061768
$ws.Range("J2").Value = 178.1206817626953
$ws.Range("K2").Value = 194.7444610595703
